# Update the date heading (first paragraph)
$d = $word.ActiveDocument
$d.Paragraphs.Item(1).Range.Text = "2023-10-29 Sunday"

# Update the multiplication answer table cells
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "23×52=1196"
$t.Cell(1, 2).Range.Text = "51×65=3315"
$t.Cell(1, 3).Range.Text = "21×50=1050"
$t.Cell(1, 4).Range.Text = "36×44=1584"
$t.Cell(1, 5).Range.Text = "40×28=1120"

$t.Cell(5, 1).Range.Text = "20×71=1420"
$t.Cell(5, 2).Range.Text = "23×76=1748"
$t.Cell(5, 3).Range.Text = "44×49=2156"
$t.Cell(5, 4).Range.Text = "44×75=3300"
$t.Cell(5, 5).Range.Text = "80×68=5440"

$t.Cell(10, 1).Range.Text = "51×45=2295"
$t.Cell(10, 2).Range.Text = "65×44=2860"
$t.Cell(10, 3).Range.Text = "87×48=4176"
$t.Cell(10, 4).Range.Text = "15×65=975"
$t.Cell(10, 5).Range.Text = "14×75=1050"

$t.Cell(15, 1).Range.Text = "98×16=1568"
$t.Cell(15, 2).Range.Text = "97×26=2522"
$t.Cell(15, 3).Range.Text = "18×39=702"
$t.Cell(15, 4).Range.Text = "56×40=2240"
$t.Cell(15, 5).Range.Text = "65×80=5200"

$t.Cell(20, 1).Range.Text = "86×73=6278"
$t.Cell(20, 2).Range.Text = "97×77=7469"
$t.Cell(20, 3).Range.Text = "64×59=3776"
$t.Cell(20, 4).Range.Text = "44×56=2464"
$t.Cell(20, 5).Range.Text = "86×64=5504"
